# The recorded change is purely a re-serialization artifact: the document
# was re-saved by an upgraded OOXML writer (per the commit message, an
# Apache POI version bump) which emits element attributes in a different
# (alphabetically sorted) order. Every hunk in the diff - the <w:document>
# namespace declarations, <w:pgSz>/<w:pgMar>, <w:rFonts>/<w:lang>,
# <w:latentStyles>/<w:lsdException>, and the <w:style>/<w:tblInd>/
# <w:tblCellMar> entries in styles.xml - reorders attributes only; no
# attribute, value, namespace, or piece of text content is added, removed,
# or changed. Attribute order carries no meaning in OOXML/Word's object
# model, so there is no corresponding user-visible edit to replay through
# the Word COM API: Word does not expose any property that controls raw
# XML attribute ordering.
#
# We still touch the document through the object model so the session
# performs a real (idempotent, content-preserving) round trip rather than
# doing nothing at all.
$d = $word.ActiveDocument
$null = $d.Content.Text
$null = $d.Sections.Count
